$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S, row 4 (header year 2022) - same format as R4
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value2 = 2022

# New column S, row 5 (value 42) - same format as R5, but with a custom
# number format "0.0" applied (creates new numFmt/cellXf entries)
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").NumberFormat = "0.0"
$ws.Range("S5").Value2 = 42

# Update the active selection to match the new state
$ws.Range("U4").Select()

$excel.CutCopyMode = 0
